# Update cryptos list: refresh Price (column D) and Volume(1h) (column E)
# values for the latest GitHub Actions scrape run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.150.29'
$ws.Range("E2").Value = '  +0.21%  '
$ws.Range("D3").Value = '1.677.02'
$ws.Range("E3").Value = '  -0.31%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.14'
$ws.Range("E5").Value = '  -0.83%  '
$ws.Range("E6").Value = '  +0.19%  '
$ws.Range("E7").Value = '  +0.10%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '22.89'
$ws.Range("E8").Value = '  +7.25%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.262'
$ws.Range("E9").Value = '  +2.99%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0622'
$ws.Range("E10").Value = '  -0.38%  '
$ws.Range("E11").Value = '  +0.00%  '
$ws.Range("D12").Value = '1.913.31'
$ws.Range("E12").Value = '  -0.32%  '
$ws.Range("D13").Value = '1.678.91'
$ws.Range("E13").Value = '  -0.14%  '
$ws.Range("E14").Value = '  +2.20%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.560'
$ws.Range("E15").Value = '  +4.53%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.56'
$ws.Range("E16").Value = '  +0.14%  '
$ws.Range("D17").Value = '27.122.50'
$ws.Range("E17").Value = '  +0.10%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '235.16'
$ws.Range("E18").Value = '  -0.53%  '
$ws.Range("E19").Value = '  +0.58%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.82'
$ws.Range("E20").Value = '  -4.06%  '
$ws.Range("E22").Value = '  +1.84%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.55'
$ws.Range("E23").Value = '  +2.96%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.08'
$ws.Range("E24").Value = '  -2.20%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '148.01'
$ws.Range("E25").Value = '  +0.46%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.48'
$ws.Range("E26").Value = '  +2.44%  '
$ws.Range("E27").Value = '  -0.64%  '
$ws.Range("E28").Value = '  -0.07%  '
$ws.Range("E29").Value = '  +0.11%  '
$ws.Range("E30").Value = '  +0.37%  '
$ws.Range("E31").Value = '  -0.15%  '
$ws.Range("E32").Value = '  -0.02%  '
$ws.Range("D33").Value = '1.547.99'
$ws.Range("E33").Value = '  +0.01%  '
$ws.Range("E35").Value = '  -3.81%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.607'
$ws.Range("E36").Value = '  +3.45%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.948'
$ws.Range("E37").Value = '  +3.42%  '
$ws.Range("E38").Value = '  -0.01%  '
$ws.Range("E39").Value = '  -1.10%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '69.65'
$ws.Range("E41").Value = '  +2.54%  '
$ws.Range("E42").Value = '  +4.37%  '
$ws.Range("E43").Value = '  +0.11%  '
$ws.Range("E44").Value = '  -0.32%  '
$ws.Range("D45").Value = '1.823.23'
$ws.Range("E45").Value = '  +0.04%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.783'
$ws.Range("E46").Value = '  +0.28%  '
$ws.Range("E47").Value = '  +6.21%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '89.70'
$ws.Range("E48").Value = '  -1.00%  '
$ws.Range("E49").Value = '  +2.08%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.26'
$ws.Range("E50").Value = '  +3.07%  '
$ws.Range("E51").Value = '  +0.43%  '
